$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.3227736666666667
$ws.Cells.Item(2, 8).Value = 0.968321
$ws.Cells.Item(2, 9).Value = 0.1416094457286952
$ws.Cells.Item(2, 10).Value = 0.1416094457286952
$ws.Cells.Item(2, 13).Value = 2.655411666666667
$ws.Cells.Item(2, 14).Value = 7.966235
$ws.Cells.Item(2, 15).Value = 0.1255433399118981
$ws.Cells.Item(2, 16).Value = 0.1255433399118982
$ws.Cells.Item(2, 17).Value = 0.8570969601594445
$ws.Cells.Item(2, 18).Value = 7.713872641435001
$ws.Cells.Item(2, 19).Value = 0.01777812277985307
$ws.Cells.Item(2, 20).Value = 0.01777812277985307

$ws.Cells.Item(3, 7).Value = 0.3227736666666667
$ws.Cells.Item(3, 8).Value = 0.968321
$ws.Cells.Item(3, 9).Value = 0.1416094457286952
$ws.Cells.Item(3, 10).Value = 0.1416094457286952
$ws.Cells.Item(3, 15).Value = 0.3702382146908386
$ws.Cells.Item(3, 16).Value = 0.3702382146908386
$ws.Cells.Item(3, 17).Value = 2.527653387022111
$ws.Cells.Item(3, 18).Value = 22.748880483199
$ws.Cells.Item(3, 19).Value = 0.0524292283699513
$ws.Cells.Item(3, 20).Value = 0.0524292283699513

$ws.Cells.Item(4, 7).Value = 0.3227736666666667
$ws.Cells.Item(4, 8).Value = 0.968321
$ws.Cells.Item(4, 9).Value = 0.1416094457286952
$ws.Cells.Item(4, 10).Value = 0.1416094457286952
$ws.Cells.Item(4, 13).Value = 6.127532
$ws.Cells.Item(4, 14).Value = 18.382596
$ws.Cells.Item(4, 15).Value = 0.2896992742608144
$ws.Cells.Item(4, 16).Value = 0.2896992742608145
$ws.Cells.Item(4, 17).Value = 1.977805971257333
$ws.Cells.Item(4, 18).Value = 17.800253741316
$ws.Cells.Item(4, 19).Value = 0.04102415365607918
$ws.Cells.Item(4, 20).Value = 0.04102415365607918

$ws.Cells.Item(5, 7).Value = 0.3227736666666667
$ws.Cells.Item(5, 8).Value = 0.968321
$ws.Cells.Item(5, 9).Value = 0.1416094457286952
$ws.Cells.Item(5, 10).Value = 0.1416094457286952
$ws.Cells.Item(5, 13).Value = 2.108791333333333
$ws.Cells.Item(5, 14).Value = 6.326373999999999
$ws.Cells.Item(5, 15).Value = 0.09970006175963861
$ws.Cells.Item(5, 16).Value = 0.09970006175963862
$ws.Cells.Item(5, 17).Value = 0.6806623108948888
$ws.Cells.Item(5, 18).Value = 6.125960798054
$ws.Cells.Item(5, 19).Value = 0.0141184704848991
$ws.Cells.Item(5, 20).Value = 0.0141184704848991

$ws.Cells.Item(6, 7).Value = 0.3227736666666667
$ws.Cells.Item(6, 8).Value = 0.968321
$ws.Cells.Item(6, 9).Value = 0.1416094457286952
$ws.Cells.Item(6, 10).Value = 0.1416094457286952
$ws.Cells.Item(6, 13).Value = 2.428579666666666
$ws.Cells.Item(6, 14).Value = 7.285739
$ws.Cells.Item(6, 15).Value = 0.1148191093768101
$ws.Cells.Item(6, 16).Value = 0.1148191093768101
$ws.Cells.Item(6, 17).Value = 0.783881563802111
$ws.Cells.Item(6, 18).Value = 7.054934074218999
$ws.Cells.Item(6, 19).Value = 0.0162594704379125
$ws.Cells.Item(6, 20).Value = 0.01625947043791251

$ws.Cells.Item(7, 9).Value = 0.8226066833587575
$ws.Cells.Item(7, 10).Value = 0.8226066833587576
$ws.Cells.Item(7, 13).Value = 2.655411666666667
$ws.Cells.Item(7, 14).Value = 7.966235
$ws.Cells.Item(7, 15).Value = 0.1255433399118981
$ws.Cells.Item(7, 16).Value = 0.1255433399118982
$ws.Cells.Item(7, 17).Value = 4.978860584373889
$ws.Cells.Item(7, 18).Value = 44.809745259365
$ws.Cells.Item(7, 19).Value = 0.1032727904627077
$ws.Cells.Item(7, 20).Value = 0.1032727904627077

$ws.Cells.Item(8, 9).Value = 0.8226066833587575
$ws.Cells.Item(8, 10).Value = 0.8226066833587576
$ws.Cells.Item(8, 15).Value = 0.3702382146908386
$ws.Cells.Item(8, 16).Value = 0.3702382146908386
$ws.Cells.Item(8, 19).Value = 0.3045604298394984
$ws.Cells.Item(8, 20).Value = 0.3045604298394984

$ws.Cells.Item(9, 9).Value = 0.8226066833587575
$ws.Cells.Item(9, 10).Value = 0.8226066833587576
$ws.Cells.Item(9, 13).Value = 6.127532
$ws.Cells.Item(9, 14).Value = 18.382596
$ws.Cells.Item(9, 15).Value = 0.2896992742608144
$ws.Cells.Item(9, 16).Value = 0.2896992742608145
$ws.Cells.Item(9, 17).Value = 11.48903875706267
$ws.Cells.Item(9, 18).Value = 103.401348813564
$ws.Cells.Item(9, 19).Value = 0.2383085591711276
$ws.Cells.Item(9, 20).Value = 0.2383085591711277

$ws.Cells.Item(10, 9).Value = 0.8226066833587575
$ws.Cells.Item(10, 10).Value = 0.8226066833587576
$ws.Cells.Item(10, 13).Value = 2.108791333333333
$ws.Cells.Item(10, 14).Value = 6.326373999999999
$ws.Cells.Item(10, 15).Value = 0.09970006175963861
$ws.Cells.Item(10, 16).Value = 0.09970006175963862
$ws.Cells.Item(10, 17).Value = 3.953954929851777
$ws.Cells.Item(10, 18).Value = 35.585594368666
$ws.Cells.Item(10, 19).Value = 0.0820139371347596
$ws.Cells.Item(10, 20).Value = 0.08201393713475963

$ws.Cells.Item(11, 9).Value = 0.8226066833587575
$ws.Cells.Item(11, 10).Value = 0.8226066833587576
$ws.Cells.Item(11, 13).Value = 2.428579666666666
$ws.Cells.Item(11, 14).Value = 7.285739
$ws.Cells.Item(11, 15).Value = 0.1148191093768101
$ws.Cells.Item(11, 16).Value = 0.1148191093768101
$ws.Cells.Item(11, 17).Value = 4.553553684411222
$ws.Cells.Item(11, 18).Value = 40.981983159701
$ws.Cells.Item(11, 19).Value = 0.09445096675066417
$ws.Cells.Item(11, 20).Value = 0.09445096675066419

$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.081563
$ws.Cells.Item(12, 8).Value = 0.244689
$ws.Cells.Item(12, 9).Value = 0.03578387091254728
$ws.Cells.Item(12, 10).Value = 0.03578387091254728
$ws.Cells.Item(12, 13).Value = 2.655411666666667
$ws.Cells.Item(12, 14).Value = 7.966235
$ws.Cells.Item(12, 15).Value = 0.1255433399118981
$ws.Cells.Item(12, 16).Value = 0.1255433399118982
$ws.Cells.Item(12, 17).Value = 0.2165833417683333
$ws.Cells.Item(12, 18).Value = 1.949250075915
$ws.Cells.Item(12, 19).Value = 0.004492426669337407
$ws.Cells.Item(12, 20).Value = 0.004492426669337408

$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.081563
$ws.Cells.Item(13, 8).Value = 0.244689
$ws.Cells.Item(13, 9).Value = 0.03578387091254728
$ws.Cells.Item(13, 10).Value = 0.03578387091254728
$ws.Cells.Item(13, 15).Value = 0.3702382146908386
$ws.Cells.Item(13, 16).Value = 0.3702382146908386
$ws.Cells.Item(13, 17).Value = 0.6387230883323333
$ws.Cells.Item(13, 18).Value = 5.748507794991
$ws.Cells.Item(13, 19).Value = 0.01324855648138893
$ws.Cells.Item(13, 20).Value = 0.01324855648138893

$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 0.081563
$ws.Cells.Item(14, 8).Value = 0.244689
$ws.Cells.Item(14, 9).Value = 0.03578387091254728
$ws.Cells.Item(14, 10).Value = 0.03578387091254728
$ws.Cells.Item(14, 13).Value = 6.127532
$ws.Cells.Item(14, 14).Value = 18.382596
$ws.Cells.Item(14, 15).Value = 0.2896992742608144
$ws.Cells.Item(14, 16).Value = 0.2896992742608145
$ws.Cells.Item(14, 17).Value = 0.499779892516
$ws.Cells.Item(14, 18).Value = 4.498019032644
$ws.Cells.Item(14, 19).Value = 0.01036656143360761
$ws.Cells.Item(14, 20).Value = 0.01036656143360761

$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 0.081563
$ws.Cells.Item(15, 8).Value = 0.244689
$ws.Cells.Item(15, 9).Value = 0.03578387091254728
$ws.Cells.Item(15, 10).Value = 0.03578387091254728
$ws.Cells.Item(15, 13).Value = 2.108791333333333
$ws.Cells.Item(15, 14).Value = 6.326373999999999
$ws.Cells.Item(15, 15).Value = 0.09970006175963861
$ws.Cells.Item(15, 16).Value = 0.09970006175963862
$ws.Cells.Item(15, 17).Value = 0.1719993475206666
$ws.Cells.Item(15, 18).Value = 1.547994127686
$ws.Cells.Item(15, 19).Value = 0.003567654139979899
$ws.Cells.Item(15, 20).Value = 0.0035676541399799

$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 0.081563
$ws.Cells.Item(16, 8).Value = 0.244689
$ws.Cells.Item(16, 9).Value = 0.03578387091254728
$ws.Cells.Item(16, 10).Value = 0.03578387091254728
$ws.Cells.Item(16, 13).Value = 2.428579666666666
$ws.Cells.Item(16, 14).Value = 7.285739
$ws.Cells.Item(16, 15).Value = 0.1148191093768101
$ws.Cells.Item(16, 16).Value = 0.1148191093768101
$ws.Cells.Item(16, 17).Value = 0.2165833417683333
$ws.Cells.Item(16, 18).Value = 1.949250075915
$ws.Cells.Item(16, 19).Value = 0.004492426669337407
$ws.Cells.Item(16, 20).Value = 0.004492426669337408
